$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# Row 71: Mineros de Zacatecas vs Cancún - Draw, Fallo
$ws.Range("L71").Value = "Completed"
$ws.Range("M71").Value = "Draw"
$ws.Range("N71").Value = "Fallo"
$ws.Range("O71").Value = -1.3
$ws.Range("P71").Value = -100
$ws.Range("Q71").Value = "2025-09-07 04:25:41"

# Row 72: Dorados vs Irapuato - Away Win, Acierto
$ws.Range("L72").Value = "Completed"
$ws.Range("M72").Value = "Away Win"
$ws.Range("N72").Value = "Acierto"
$ws.Range("O72").Value = 0.6
$ws.Range("P72").Value = 120
$ws.Range("Q72").Value = "2025-09-07 04:25:41"

# Row 73: Tlaxcala vs Alebrijes de Oaxaca - Home Win, Acierto
$ws.Range("L73").Value = "Completed"
$ws.Range("M73").Value = "Home Win"
$ws.Range("N73").Value = "Acierto"
$ws.Range("O73").Value = 1.58
$ws.Range("P73").Value = 75
$ws.Range("Q73").Value = "2025-09-07 04:25:41"
